# Apply the commit's data changes to column F ("dSF") for the rows
# whose values were "repulled" in this update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = 15
    9  = -2
    14 = -2
    15 = 4
    16 = 4
    18 = -1
    19 = -5
    20 = 1
    22 = -4
    23 = -4
    24 = -2
    26 = 0
    28 = -2
    29 = 3
    30 = 1
    35 = -2
    37 = -3
    39 = 1
    44 = -8
    45 = -6
    50 = -2
    56 = -6
    57 = 3
    58 = 3
    59 = 8
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
